$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.949.53'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '2.779.12'
$ws.Range('E3').Value = '  -1.61%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '356.00'
$ws.Range('E5').Value = '  +0.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '108.92'
$ws.Range('E6').Value = '  -4.30%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.561'
$ws.Range('E7').Value = '  +1.98%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.589'
$ws.Range('E9').Value = '  -1.97%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.12'
$ws.Range('E10').Value = '  -4.49%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0850'
$ws.Range('E11').Value = '  -0.10%  '
$ws.Range('E12').Value = '  +0.80%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.35'
$ws.Range('E13').Value = '  -3.68%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.59'
$ws.Range('E14').Value = '  -1.58%  '
$ws.Range('D15').Value = '3.210.71'
$ws.Range('E15').Value = '  -1.17%  '
$ws.Range('D16').Value = '2.790.70'
$ws.Range('E16').Value = '  -1.10%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.933'
$ws.Range('E17').Value = '  +4.01%  '
$ws.Range('D18').Value = '51.811.81'
$ws.Range('E18').Value = '  -0.12%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.38'
$ws.Range('E19').Value = '  +0.56%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.13'
$ws.Range('E20').Value = '  -1.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.03'
$ws.Range('E21').Value = '  -3.95%  '
$ws.Range('D22').Value = '0.0₃0974'
$ws.Range('E22').Value = '  -2.62%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '274.27'
$ws.Range('E23').Value = '  +1.67%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '69.71'
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.73'
$ws.Range('E25').Value = '  -2.41%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.53'
$ws.Range('E26').Value = '  -0.87%  '
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.12'
$ws.Range('E28').Value = '  -1.95%  '
$ws.Range('E29').Value = '  -1.13%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.143'
$ws.Range('E30').Value = '  +2.00%  '
$ws.Range('B31').Value = 'VeChain'
$ws.Range('C31').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0465'
$ws.Range('E31').Value = '  +3.19%  '
$ws.Range('B32').Value = 'OKB'
$ws.Range('C32').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '51.38'
$ws.Range('E32').Value = '  +1.11%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '33.72'
$ws.Range('E33').Value = '  -0.31%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.70'
$ws.Range('E34').Value = '  -2.58%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.31'
$ws.Range('E35').Value = '  +9.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0836'
$ws.Range('E36').Value = '  +0.30%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.19'
$ws.Range('E38').Value = '  -0.49%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.05'
$ws.Range('E39').Value = '  -1.79%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.99'
$ws.Range('E40').Value = '  -4.90%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.53'
$ws.Range('E41').Value = '  -1.84%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.115'
$ws.Range('E42').Value = '  -0.52%  '
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.24'
$ws.Range('E43').Value = '  -3.02%  '
$ws.Range('B44').Value = 'Monero'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '120.91'
$ws.Range('E44').Value = '  -5.76%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.68'
$ws.Range('E45').Value = '  -8.42%  '
$ws.Range('D46').Value = '2.053.07'
$ws.Range('E46').Value = '  -1.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.24'
$ws.Range('E47').Value = '  -2.92%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.27'
$ws.Range('E48').Value = '  -1.92%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.69'
$ws.Range('E49').Value = '  -0.14%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.922'
$ws.Range('E50').Value = '  -3.90%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.88'
$ws.Range('E51').Value = '  -0.19%  '
